{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst items = paras.items;\n\n// 1. Turn the 3rd paragraph (empty) into the new class declaration line.\nitems[2].insertText(\"class TestPasswordValidation(unittest.TestCase):\", \"Replace\");\n\n// 2. Drop one of the two blank paragraphs that followed\n//    `password = \"abc1$\"` (collapse blank-blank -> blank).\nitems[6].delete();\n\n// 3. Drop one of the two blank paragraphs that followed\n//    `password = \"abcd@xyz\"` / test_no_number block.\nitems[11].delete();\n\n// 4. The blank paragraph after `password = \"abcd1234\"` becomes a 4-space\n//    indent line, followed by a brand-new test_no_number_2 block.\nitems[14].insertText(\"    \", \"Replace\");\nconst newPara1 = items[14].insertParagraph(\"    def test_no_number_2(self):\", \"After\");\nnewPara1.insertParagraph('        password = \"abcd@xyz\"', \"After\");\n\n// 5. Rename test_valid_8_chars -> test_valid_password, and drop an extra\n//    blank separator paragraph before it.\nitems[19].delete();\nitems[16].insertText(\"    def test_valid_password(self):\", \"Replace\");\n\n// 6. Rename test_valid_gt_8_chars -> test_valid_long_password.\nitems[23].delete();\nitems[20].insertText(\"    def test_valid_long_password(self):\", \"Replace\");\n\n// 7. Rename test_valid_new1 -> test_client_input_1.\nitems[27].delete();\nitems[24].insertText(\"    def test_client_input_1(self):\", \"Replace\");\n\n// 8. Rename test_invalid_new2 -> test_client_input_2.\nitems[31].delete();\nitems[28].insertText(\"    def test_client_input_2(self):\", \"Replace\");\n\n// 9. Rename test_invalid_new3 -> test_client_input_3.\nitems[35].delete();\nitems[32].insertText(\"    def test_client_input_3(self):\", \"Replace\");\n\n// 10. Rename test_valid_new4 -> test_client_input_4.\nitems[39].delete();\nitems[36].insertText(\"    def test_client_input_4(self):\", \"Replace\");\n\n// 11. Collapse the long comment + manual validation body of\n//     is_valid_password() into a single comment line.\nitems[41].insertText(\"    # Password validation logic\", \"Replace\");\nitems[42].delete();\nitems[43].delete();\nitems[44].delete();\nitems[45].delete();\nitems[46].delete();\nitems[47].delete();\nitems[48].delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Word COM paragraph/range positions drift as earlier edits shift the\n# character stream, so every paragraph below is looked up FRESH, by its\n# original (1-based) index, at the moment it is needed. Processing strictly\n# from the bottom of the document upward guarantees that every remaining\n# lookup still sits above all the edits already made, so its index is still\n# correct when we fetch it.\n\n# 11. Collapse the long comment + manual validation body of\n#     is_valid_password() into a single comment line (originally\n#     paragraphs 42-49: comment, comment, if, return, if, return, if,\n#     return).\n$d.Paragraphs.Item(49).Range.Delete()\n$d.Paragraphs.Item(48).Range.Delete()\n$d.Paragraphs.Item(47).Range.Delete()\n$d.Paragraphs.Item(46).Range.Delete()\n$d.Paragraphs.Item(45).Range.Delete()\n$d.Paragraphs.Item(44).Range.Delete()\n$d.Paragraphs.Item(43).Range.Delete()\n$d.Paragraphs.Item(42).Range.Text = \"    # Password validation logic\"\n\n# 10. Rename test_valid_new4 -> test_client_input_4, dropping one of the\n#     two blank separators before it.\n$d.Paragraphs.Item(40).Range.Delete()\n$d.Paragraphs.Item(37).Range.Text = \"    def test_client_input_4(self):\"\n\n# 9. Rename test_invalid_new3 -> test_client_input_3.\n$d.Paragraphs.Item(36).Range.Delete()\n$d.Paragraphs.Item(33).Range.Text = \"    def test_client_input_3(self):\"\n\n# 8. Rename test_invalid_new2 -> test_client_input_2.\n$d.Paragraphs.Item(32).Range.Delete()\n$d.Paragraphs.Item(29).Range.Text = \"    def test_client_input_2(self):\"\n\n# 7. Rename test_valid_new1 -> test_client_input_1.\n$d.Paragraphs.Item(28).Range.Delete()\n$d.Paragraphs.Item(25).Range.Text = \"    def test_client_input_1(self):\"\n\n# 6. Rename test_valid_gt_8_chars -> test_valid_long_password.\n$d.Paragraphs.Item(24).Range.Delete()\n$d.Paragraphs.Item(21).Range.Text = \"    def test_valid_long_password(self):\"\n\n# 5. Rename test_valid_8_chars -> test_valid_password.\n$d.Paragraphs.Item(20).Range.Delete()\n$d.Paragraphs.Item(17).Range.Text = \"    def test_valid_password(self):\"\n\n# 4. The blank paragraph after `password = \"abcd1234\"` becomes a 4-space\n#    indent line, followed by a brand-new test_no_number_2 block.\n$pIndent = $d.Paragraphs.Item(15)\n$pIndent.Range.Text = \"    \"\n$pIndent.Range.InsertParagraphAfter()\n$pDef = $pIndent.Next()\n$pDef.Range.Text = \"    def test_no_number_2(self):\"\n$pDef.Range.InsertParagraphAfter()\n$pPwd = $pDef.Next()\n$pPwd.Range.Text = '        password = \"abcd@xyz\"'\n\n# 3. Drop one of the two blank paragraphs after the test_no_number block.\n$d.Paragraphs.Item(12).Range.Delete()\n\n# 2. Drop one of the two blank paragraphs after `password = \"abc1$\"`.\n$d.Paragraphs.Item(7).Range.Delete()\n\n# 1. Turn the blank paragraph into the new class declaration line.\n$d.Paragraphs.Item(3).Range.Text = \"class TestPasswordValidation(unittest.TestCase):\"\n"}
